# Add save load script
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Todo ")

# Update existing "fish" creator entries on rows 30-31 to "Fish"
$ws.Range("D30").Value = "Fish"
$ws.Range("D31").Value = "Fish"

# Add a new row 32 for the "Refactor code" todo item
$ws.Range("A32").Value = "Refactor code"
$ws.Range("B32").Value = "Refactor the code not done by ourselves"
$ws.Range("D32").Value = "Fish"
$ws.Range("F32").Value = "18 Jan"

$ws.Range("G32").Select()
